$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy H1's format (bold, border, centered) onto I1:J1, then set values
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: I column, J column
$data = @(
    @(8, 8),   # row 2
    @(7, 8),   # row 3
    @(9, 9),   # row 4
    @(9, 9),   # row 5
    @(9, 9),   # row 6
    @(6, 7),   # row 7
    @(5, 5),   # row 8
    @(9, 9),   # row 9
    @(5, 6),   # row 10
    @(5, 6),   # row 11
    @(3, 4),   # row 12
    @(8, 8),   # row 13
    @(8, 8),   # row 14
    @(9, 9),   # row 15
    @(7, 9),   # row 16
    @(5, 5)    # row 17
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
